# Add two new columns, "I0" (I) and "IF" (J), matching the existing
# header style used by the other header cells (copy format from H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy H1's formatting onto the two new header cells
# then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (no explicit style - plain values)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
